# Weapon spell slots & gem sockets: replace the old "type + 8 boolean
# modifier flags" strings with a simpler "type, max modifiers" scheme, and
# refresh the slot-1 description string. Also update view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")

# Row 4 - it_we_icestaff: slots now store "type, max modifiers".
$ws.Range("AC4").Value = "1, 3"
$ws.Range("AD4").Value = "1, 3"
$ws.Range("AE4").Value = "4, 3"
$ws.Range("AF4").Value = "4, 3"
$ws.Range("AG4").Value = "3, 3"

# Row 5 - it_we_rustysword: same update, one fewer slot used.
$ws.Range("AC5").Value = "1, 3"
$ws.Range("AD5").Value = "1, 3"
$ws.Range("AE5").Value = "3, 3"
$ws.Range("AF5").Value = "2, 3"

# Header description for slot 1 (AC1) - translation update.
$ws.Range("AC1").Value = "slot 1 (the first entry is the type (elemental, twilight, necromancy, divine, illusion) and after this, the number shows the maximal modifiers for this slot"

# Restore the saved view/selection state recorded in the workbook.
$ws.Application.ActiveWindow.ScrollColumn = 15
$ws.Range("AG8").Select()
